$d = $word.ActiveDocument

# 1. {Lección:Haz una pausa} -> {Sesión:Haz una pausa}
$d.Content.Find.Execute(
    "{Lección:Haz una pausa}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{Sesión:Haz una pausa}", 2) | Out-Null

# 2. ¡Hola! ¿Cómo te sientes ahora?  -> ¡Hola! ¿Cómo te sientes en este momento?
$d.Content.Find.Execute(
    "¡Hola! ¿Cómo te sientes ahora? ", $true, $false, $false, $false, $false,
    $true, 1, $false, "¡Hola! ¿Cómo te sientes en este momento? ", 2) | Out-Null

# 3. Antes de empezar con el Crianza con Conciencia+, hagamos juntos una breve pausa.
#    -> Hagamos una breve pausa juntos, antes de comenzar Crianza con Conciencia+.
$d.Content.Find.Execute(
    "Antes de empezar con el Crianza con Conciencia+, hagamos juntos una breve pausa. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Hagamos una breve pausa juntos, antes de comenzar Crianza con Conciencia+. ", 2) | Out-Null

# 4. Haga una pausa -> Haz una pausa
$d.Content.Find.Execute(
    "Haga una pausa", $true, $false, $false, $false, $false,
    $true, 1, $false, "Haz una pausa", 2) | Out-Null

# 5. Respira hondo.  -> Respira profundo.
$d.Content.Find.Execute(
    "Respira hondo. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Respira profundo. ", 2) | Out-Null

# 6. The SECOND occurrence of "[pausa]" (the one between "En;" and "y fuera; ")
#    becomes "[pause]". All other "[pausa]"/"[pausa] " instances stay untouched,
#    so we scan match-by-match instead of a document-wide replace-all.
$rng = $d.Content
$rng.Start = 0
$matchIndex = 0
$targetMatch = 2
$replaced = $false
while (-not $replaced) {
    $found = $rng.Find.Execute(
        "[pausa]", $true, $false, $false, $false, $false,
        $true, 0, $false, "", 0)
    if (-not $found) { break }
    $matchIndex = $matchIndex + 1
    if ($matchIndex -eq $targetMatch) {
        $rng.Text = "[pause]"
        $replaced = $true
    } else {
        $rng.Collapse(0)
    }
}
